$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update fluid_mass (B11) with the new value; Q_cool (B12) recalculates
# automatically via its formula fluid_mass*c_water*dt_cool.
$ws.Range("B11").Value = 1.4137154999999999

$excel.CalculateFullRebuild()
